# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G (header "K") values were recalculated; update rows 2-55 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(1,0,3,0,3,0,5,2,2,2,1,2,4,1,4,0,0,1,3,3,2,1,2,1,2,2,3,1,1,3,8,2,2,0,2,4,1,0,1,2,2,2,5,1,2,2,2,1,0,1,0,2,2,2)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
